$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values that look like plain decimals get auto-typed as
# numbers by Excel on assignment; the source data keeps them as text, so we
# force a text format, assign, then drop back to the Normal style so no
# stray number-format style sticks to the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.332.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.966.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.62%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.541"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.89%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.24%  "

$ws.Range("E11").Value = "  +0.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0842"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.428.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.958.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.969"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.287.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +22.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.171"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.71%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.114"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.01%  "

$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.92%  "

$ws.Range("E34").Value = "  -1.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0446"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.08%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.54%  "

$ws.Range("E42").Value = "  +3.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.293"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +26.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.83%  "

$ws.Range("E47").Value = "  +3.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.043.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0350"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.90%  "

$ws.Range("E51").Value = "  +2.90%  "

